# Rebuild the "Periodo Mora" detail table (rows 16-29).
# The data is restructured from "grouped by period" (each period listing both
# workers) to "grouped by worker" (each worker listing all of their periods,
# most recent period first). JOSE NEDER HERNANDEZ PEÑA's records now come
# first (rows 16-22), followed by JOSE HERMES HERNANDEZ OTALORA (rows 23-29).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("CC", "1007270769", "JOSE NEDER HERNANDEZ PEÑA",      "2406", 18560, 1160000),
    @("CC", "1007270769", "JOSE NEDER HERNANDEZ PEÑA",      "2405", 46400, 1160000),
    @("CC", "1007270769", "JOSE NEDER HERNANDEZ PEÑA",      "2404", 46400, 1160000),
    @("CC", "1007270769", "JOSE NEDER HERNANDEZ PEÑA",      "2403", 46400, 1160000),
    @("CC", "1007270769", "JOSE NEDER HERNANDEZ PEÑA",      "2402", 46400, 1160000),
    @("CC", "1007270769", "JOSE NEDER HERNANDEZ PEÑA",      "2401", 46400, 1160000),
    @("CC", "1007270769", "JOSE NEDER HERNANDEZ PEÑA",      "2312", 46400, 1160000),
    @("CC", "1007270040", "JOSE HERMES HERNANDEZ OTALORA",  "2406", 37120, 2320000),
    @("CC", "1007270040", "JOSE HERMES HERNANDEZ OTALORA",  "2405", 92800, 2320000),
    @("CC", "1007270040", "JOSE HERMES HERNANDEZ OTALORA",  "2404", 92800, 2320000),
    @("CC", "1007270040", "JOSE HERMES HERNANDEZ OTALORA",  "2403", 92800, 2320000),
    @("CC", "1007270040", "JOSE HERMES HERNANDEZ OTALORA",  "2402", 92800, 2320000),
    @("CC", "1007270040", "JOSE HERMES HERNANDEZ OTALORA",  "2401", 92800, 2320000),
    @("CC", "1007270040", "JOSE HERMES HERNANDEZ OTALORA",  "2312", 92800, 2320000)
)

$startRow = 16
for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $startRow + $i
    $vals = $rows[$i]
    $ws.Cells.Item($r, 2).Value = $vals[0]
    $ws.Cells.Item($r, 3).Value = $vals[1]
    $ws.Cells.Item($r, 4).Value = $vals[2]
    $ws.Cells.Item($r, 5).Value = $vals[3]
    $ws.Cells.Item($r, 6).Value = $vals[4]
    $ws.Cells.Item($r, 7).Value = $vals[5]
}
